$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 existing values
$ws.Range("B2").Value = 22
$ws.Range("C2").Value = 29

# New inventory rows 3-10 (SKU, Sold, Price)
$data = @(
    @{ Row = 3;  Sku = "SKU102"; Sold = 10; Price = 20 },
    @{ Row = 4;  Sku = "SKU103"; Sold = 35; Price = 49 },
    @{ Row = 5;  Sku = "SKU104"; Sold = 20; Price = 155 },
    @{ Row = 6;  Sku = "SKU105"; Sold = 20; Price = 15 },
    @{ Row = 7;  Sku = "SKU106"; Sold = 25; Price = 77 },
    @{ Row = 8;  Sku = "SKU107"; Sold = 33; Price = 440 },
    @{ Row = 9;  Sku = "SKU108"; Sold = 35; Price = 278 },
    @{ Row = 10; Sku = "SKU109"; Sold = 20; Price = 43 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Sku
    $ws.Cells.Item($r, 2).Value = $item.Sold
    $ws.Cells.Item($r, 3).Value = $item.Price
}

# Move active selection to C9, matching the last user-selected cell
$ws.Range("C9").Select()
